$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 302
$ws.Range("B302").Value = 6943423
$ws.Range("F302").Value = "AEK Larnaca"
$ws.Range("G302").Value = "Omonia Nicosia"
$ws.Range("H302").Value = 2
$ws.Range("I302").Value = 1
$ws.Range("J302").Value = "H"
$ws.Range("K302").Value = 1.7
$ws.Range("L302").Value = 3.75
$ws.Range("M302").Value = 4.75
$ws.Range("N302").Value = 1.75
$ws.Range("O302").Value = 3.6
$ws.Range("P302").Value = 5
$ws.Range("Q302").Value = -0.75
$ws.Range("R302").Value = 1.925
$ws.Range("S302").Value = 1.875
$ws.Range("T302").Value = 2.5
$ws.Range("U302").Value = 1.925
$ws.Range("V302").Value = 1.875
$ws.Range("W302").Value = 0.75
$ws.Range("X302").Value = -1
$ws.Range("Y302").Value = -1
$ws.Range("Z302").Value = 0.4625
$ws.Range("AA302").Value = -0.5
$ws.Range("AB302").Value = 0.925
$ws.Range("AC302").Value = -1

# Row 303
$ws.Range("B303").Value = 6942205
$ws.Range("F303").Value = "Pafos FC"
$ws.Range("G303").Value = "Apollon Limassol"
$ws.Range("H303").Value = 1
$ws.Range("I303").Value = 1
$ws.Range("J303").Value = "D"
$ws.Range("K303").Value = 1.85
$ws.Range("L303").Value = 3.5
$ws.Range("M303").Value = 4.2
$ws.Range("N303").Value = 1.85
$ws.Range("O303").Value = 3.5
$ws.Range("P303").Value = 4.2
$ws.Range("Q303").Value = -0.5
$ws.Range("R303").Value = 1.9
$ws.Range("S303").Value = 1.9
$ws.Range("T303").Value = 2.25
$ws.Range("U303").Value = 1.85
$ws.Range("V303").Value = 1.95
$ws.Range("W303").Value = -1
$ws.Range("X303").Value = 2.5
$ws.Range("Y303").Value = -1
$ws.Range("Z303").Value = -1
$ws.Range("AA303").Value = 0.8999999999999999
$ws.Range("AB303").Value = -0.5
$ws.Range("AC303").Value = 0.475

# Row 304
$ws.Range("B304").Value = 6943426
$ws.Range("F304").Value = "Aris Limassol"
$ws.Range("G304").Value = "Othellos Athienou"
$ws.Range("H304").Value = 4
$ws.Range("I304").Value = 0
$ws.Range("J304").Value = "H"
$ws.Range("K304").Value = 1.222
$ws.Range("L304").Value = 6
$ws.Range("M304").Value = 12
$ws.Range("N304").Value = 1.181
$ws.Range("O304").Value = 7
$ws.Range("P304").Value = 17
$ws.Range("Q304").Value = -2
$ws.Range("R304").Value = 1.775
$ws.Range("S304").Value = 2.025
$ws.Range("T304").Value = 3.25
$ws.Range("U304").Value = 1.825
$ws.Range("V304").Value = 1.975
$ws.Range("W304").Value = 0.181
$ws.Range("X304").Value = -1
$ws.Range("Y304").Value = -1
$ws.Range("Z304").Value = 0.7749999999999999
$ws.Range("AA304").Value = -1
$ws.Range("AB304").Value = 0.825
$ws.Range("AC304").Value = -1

# Row 306
$ws.Range("B306").Value = 6943427
$ws.Range("F306").Value = "AEL Limassol"
$ws.Range("G306").Value = "Nea Salamis Famagusta"
$ws.Range("H306").Value = 3
$ws.Range("I306").Value = 1
$ws.Range("J306").Value = "H"
$ws.Range("K306").Value = 2.5
$ws.Range("L306").Value = 3.4
$ws.Range("M306").Value = 2.7
$ws.Range("N306").Value = 2.15
$ws.Range("O306").Value = 3.6
$ws.Range("P306").Value = 3.2
$ws.Range("Q306").Value = -0.25
$ws.Range("R306").Value = 1.9
$ws.Range("S306").Value = 1.9
$ws.Range("T306").Value = 2.75
$ws.Range("U306").Value = 1.975
$ws.Range("V306").Value = 1.825
$ws.Range("W306").Value = 1.15
$ws.Range("X306").Value = -1
$ws.Range("Y306").Value = -1
$ws.Range("Z306").Value = 0.8999999999999999
$ws.Range("AA306").Value = -1
$ws.Range("AB306").Value = 0.9750000000000001
$ws.Range("AC306").Value = -1

# Row 307
$ws.Range("B307").Value = 6943425
$ws.Range("F307").Value = "Doxa Katokopias"
$ws.Range("G307").Value = "Ethnikos Achnas"
$ws.Range("H307").Value = 0
$ws.Range("I307").Value = 2
$ws.Range("J307").Value = "A"
$ws.Range("K307").Value = 3
$ws.Range("L307").Value = 3.4
$ws.Range("M307").Value = 2.3
$ws.Range("N307").Value = 2.9
$ws.Range("O307").Value = 3.4
$ws.Range("P307").Value = 2.4
$ws.Range("Q307").Value = 0.25
$ws.Range("R307").Value = 1.75
$ws.Range("S307").Value = 2.05
$ws.Range("T307").Value = 2.25
$ws.Range("U307").Value = 1.8
$ws.Range("V307").Value = 2
$ws.Range("W307").Value = -1
$ws.Range("X307").Value = -1
$ws.Range("Y307").Value = 1.4
$ws.Range("Z307").Value = -1
$ws.Range("AA307").Value = 1.05
$ws.Range("AB307").Value = -0.5
$ws.Range("AC307").Value = 0.5

# Row 310
$ws.Range("B310").Value = 6943430
$ws.Range("F310").Value = "Othellos Athienou"
$ws.Range("G310").Value = "Doxa Katokopias"
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = 0
$ws.Range("J310").Value = "D"
$ws.Range("K310").Value = 2.2
$ws.Range("L310").Value = 3.4
$ws.Range("M310").Value = 3.2
$ws.Range("N310").Value = 2.2
$ws.Range("O310").Value = 3.4
$ws.Range("P310").Value = 3.1
$ws.Range("Q310").Value = -0.25
$ws.Range("R310").Value = 1.9
$ws.Range("S310").Value = 1.9
$ws.Range("T310").Value = 2.5
$ws.Range("U310").Value = 1.9
$ws.Range("V310").Value = 1.9
$ws.Range("W310").Value = -1
$ws.Range("X310").Value = 2.4
$ws.Range("Y310").Value = -1
$ws.Range("Z310").Value = -0.5
$ws.Range("AA310").Value = 0.45
$ws.Range("AB310").Value = -1
$ws.Range("AC310").Value = 0.8999999999999999

# Row 311
$ws.Range("B311").Value = 6943431
$ws.Range("F311").Value = "Ethnikos Achnas"
$ws.Range("G311").Value = "AE Zakakiou"
$ws.Range("H311").Value = 4
$ws.Range("I311").Value = 1
$ws.Range("J311").Value = "H"
$ws.Range("K311").Value = 1.85
$ws.Range("L311").Value = 3.6
$ws.Range("M311").Value = 4
$ws.Range("N311").Value = 1.909
$ws.Range("O311").Value = 3.8
$ws.Range("P311").Value = 3.8
$ws.Range("Q311").Value = -0.5
$ws.Range("R311").Value = 1.9
$ws.Range("S311").Value = 1.9
$ws.Range("T311").Value = 2.75
$ws.Range("U311").Value = 1.775
$ws.Range("V311").Value = 2.025
$ws.Range("W311").Value = 0.909
$ws.Range("X311").Value = -1
$ws.Range("Y311").Value = -1
$ws.Range("Z311").Value = 0.8999999999999999
$ws.Range("AA311").Value = -1
$ws.Range("AB311").Value = 0.7749999999999999
$ws.Range("AC311").Value = -1

# Row 312
$ws.Range("B312").Value = 6943434
$ws.Range("F312").Value = "Apoel Nicosia"
$ws.Range("G312").Value = "AEK Larnaca"
$ws.Range("H312").Value = 2
$ws.Range("I312").Value = 2
$ws.Range("J312").Value = "D"
$ws.Range("K312").Value = 1.909
$ws.Range("L312").Value = 3.4
$ws.Range("M312").Value = 4
$ws.Range("N312").Value = 2.05
$ws.Range("O312").Value = 3.4
$ws.Range("P312").Value = 3.8
$ws.Range("Q312").Value = -0.5
$ws.Range("R312").Value = 2.025
$ws.Range("S312").Value = 1.775
$ws.Range("T312").Value = 2.5
$ws.Range("U312").Value = 1.85
$ws.Range("V312").Value = 1.95
$ws.Range("W312").Value = -1
$ws.Range("X312").Value = 2.4
$ws.Range("Y312").Value = -1
$ws.Range("Z312").Value = -1
$ws.Range("AA312").Value = 0.7749999999999999
$ws.Range("AB312").Value = 0.8500000000000001
$ws.Range("AC312").Value = -1

# Row 313
$ws.Range("B313").Value = 6943433
$ws.Range("F313").Value = "APK Karmotissa"
$ws.Range("G313").Value = "AEL Limassol"
$ws.Range("H313").Value = 1
$ws.Range("I313").Value = 3
$ws.Range("J313").Value = "A"
$ws.Range("K313").Value = 2.8
$ws.Range("L313").Value = 3.25
$ws.Range("M313").Value = 2.5
$ws.Range("N313").Value = 3.1
$ws.Range("O313").Value = 3.4
$ws.Range("P313").Value = 2.3
$ws.Range("Q313").Value = 0.25
$ws.Range("R313").Value = 1.8
$ws.Range("S313").Value = 2
$ws.Range("T313").Value = 2.5
$ws.Range("U313").Value = 1.85
$ws.Range("V313").Value = 1.95
$ws.Range("W313").Value = -1
$ws.Range("X313").Value = -1
$ws.Range("Y313").Value = 1.3
$ws.Range("Z313").Value = -1
$ws.Range("AA313").Value = 1
$ws.Range("AB313").Value = 0.8500000000000001
$ws.Range("AC313").Value = -1

# Row 314
$ws.Range("B314").Value = 7647734
$ws.Range("F314").Value = "Nea Salamis Famagusta"
$ws.Range("G314").Value = "Aris Limassol"
$ws.Range("H314").Value = 0
$ws.Range("I314").Value = 3
$ws.Range("J314").Value = "A"
$ws.Range("K314").Value = 6
$ws.Range("L314").Value = 4.333
$ws.Range("M314").Value = 1.5
$ws.Range("N314").Value = 5
$ws.Range("O314").Value = 4
$ws.Range("P314").Value = 1.65
$ws.Range("Q314").Value = 0.75
$ws.Range("R314").Value = 2.025
$ws.Range("S314").Value = 1.775
$ws.Range("T314").Value = 2.75
$ws.Range("U314").Value = 1.8
$ws.Range("V314").Value = 2
$ws.Range("W314").Value = -1
$ws.Range("X314").Value = -1
$ws.Range("Y314").Value = 0.6499999999999999
$ws.Range("Z314").Value = -1
$ws.Range("AA314").Value = 0.7749999999999999
$ws.Range("AB314").Value = 0.4
$ws.Range("AC314").Value = -0.5
